$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Create Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Populate Sheet2 (FirstName / LastName / Mobile / DOB) ---
$ws2.Cells.Item(1,1).Value = "FirstName"
$ws2.Cells.Item(1,2).Value = "LastName"
$ws2.Cells.Item(1,3).Value = "Mobile"
$ws2.Cells.Item(1,4).Value = "DOB"

$ws2.Cells.Item(2,1).Value = "Teskyer"
$ws2.Cells.Item(2,2).Value = "sharma"
$ws2.Cells.Item(2,3).Value = 657584848
$ws2.Cells.Item(2,4).Value = 36142
$ws2.Cells.Item(2,4).NumberFormat = "mm-dd-yy"

# Stamp the same date format onto D3:D5 by copying D2's format, so every
# date cell shares a single cellXf entry instead of each write minting a
# brand-new (duplicate) style record.
$ws2.Cells.Item(2,4).Copy()
$ws2.Range("D3:D5").PasteSpecial(-4122) | Out-Null

$ws2.Cells.Item(3,1).Value = "johnson"
$ws2.Cells.Item(3,2).Value = "johnso"
$ws2.Cells.Item(3,3).Value = 4343434434
$ws2.Cells.Item(3,4).Value = 29283

$ws2.Cells.Item(4,1).Value = "Janel"
$ws2.Cells.Item(4,2).Value = "Tokeyo"
$ws2.Cells.Item(4,3).Value = 3232323234
$ws2.Cells.Item(4,4).Value = 33948

$ws2.Cells.Item(5,1).Value = "Cool"
$ws2.Cells.Item(5,2).Value = "Boy"
$ws2.Cells.Item(5,3).Value = 3232434545
$ws2.Cells.Item(5,4).Value = 34369

# Widen the Mobile column to fit its long numbers (matches the author's
# "best fit" auto-sized column).
$ws2.Columns.Item(3).ColumnWidth = 10.33

$ws2.Range("D5").Select() | Out-Null

# --- Populate Sheet1's new columns/rows (Age / Salary + more people) ---
$ws1.Activate()

$ws1.Cells.Item(1,3).Value = "Age"
$ws1.Cells.Item(1,4).Value = "Salary"

$ws1.Cells.Item(2,3).Value = 23
$ws1.Cells.Item(2,4).Value = 120000

$ws1.Cells.Item(3,1).Value = "Janel"
$ws1.Cells.Item(3,2).Value = "Thakkar"
$ws1.Cells.Item(3,3).Value = 32
$ws1.Cells.Item(3,4).Value = 90000

$ws1.Cells.Item(4,1).Value = "Prashuv"
$ws1.Cells.Item(4,2).Value = "Johnson"
$ws1.Cells.Item(4,3).Value = 19
$ws1.Cells.Item(4,4).Value = 87000

$ws1.Cells.Item(5,1).Value = "Aaditya"
$ws1.Cells.Item(5,2).Value = "Dhakal"
$ws1.Cells.Item(5,3).Value = 15
$ws1.Cells.Item(5,4).Value = 85000

$ws1.Cells.Item(6,1).Value = "Rabi"
$ws1.Cells.Item(6,2).Value = "Lamichhane"
$ws1.Cells.Item(6,3).Value = 48
$ws1.Cells.Item(6,4).Value = 100000

$ws1.Range("D6").Select() | Out-Null
